$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 16.12468433333333
$ws.Range("H2").Value = 48.374053
$ws.Range("I2").Value = 0.2955490655206278
$ws.Range("J2").Value = 0.2955490655206279
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 46.05975733333333
$ws.Range("N2").Value = 138.179272
$ws.Range("O2").Value = 0.8970651351272991
$ws.Range("P2").Value = 0.897065135127299
$ws.Range("Q2").Value = 742.699047469935
$ws.Range("R2").Value = 6684.291427229416
$ws.Range("S2").Value = 0.265126762398009
$ws.Range("T2").Value = 0.265126762398009

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 16.12468433333333
$ws.Range("H3").Value = 48.374053
$ws.Range("I3").Value = 0.2955490655206278
$ws.Range("J3").Value = 0.2955490655206279
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.770761666666667
$ws.Range("N3").Value = 8.312285
$ws.Range("O3").Value = 0.05396367312415441
$ws.Range("P3").Value = 0.0539636731241544
$ws.Range("Q3").Value = 44.67765723790055
$ws.Range("R3").Value = 402.098915141105
$ws.Range("S3").Value = 0.01594891316390445
$ws.Range("T3").Value = 0.01594891316390445

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 16.12468433333333
$ws.Range("H4").Value = 48.374053
$ws.Range("I4").Value = 0.2955490655206278
$ws.Range("J4").Value = 0.2955490655206279
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.802173666666667
$ws.Range("N4").Value = 5.406521
$ws.Range("O4").Value = 0.03509934175535083
$ws.Range("P4").Value = 0.03509934175535083
$ws.Range("Q4").Value = 29.05948148884588
$ws.Range("R4").Value = 261.535333399613
$ws.Range("S4").Value = 0.01037357765618309
$ws.Range("T4").Value = 0.01037357765618309

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.12468433333333
$ws.Range("H5").Value = 48.374053
$ws.Range("I5").Value = 0.2955490655206278
$ws.Range("J5").Value = 0.2955490655206279
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7122493333333333
$ws.Range("N5").Value = 2.136748
$ws.Range("O5").Value = 0.01387184999319569
$ws.Range("P5").Value = 0.01387184999319569
$ws.Range("Q5").Value = 11.48479566662711
$ws.Range("R5").Value = 103.363160999644
$ws.Range("S5").Value = 0.004099812302531315
$ws.Range("T5").Value = 0.004099812302531316

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 14.68975
$ws.Range("H6").Value = 44.06925
$ws.Range("I6").Value = 0.2692481784748309
$ws.Range("J6").Value = 0.2692481784748309
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 46.05975733333333
$ws.Range("N6").Value = 138.179272
$ws.Range("O6").Value = 0.8970651351272991
$ws.Range("P6").Value = 0.897065135127299
$ws.Range("Q6").Value = 676.6063202873332
$ws.Range("R6").Value = 6089.456882586
$ws.Range("S6").Value = 0.2415331536063033
$ws.Range("T6").Value = 0.2415331536063033

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 14.68975
$ws.Range("H7").Value = 44.06925
$ws.Range("I7").Value = 0.2692481784748309
$ws.Range("J7").Value = 0.2692481784748309
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.770761666666667
$ws.Range("N7").Value = 8.312285
$ws.Range("O7").Value = 0.05396367312415441
$ws.Range("P7").Value = 0.0539636731241544
$ws.Range("Q7").Value = 40.70179619291666
$ws.Range("R7").Value = 366.31616573625
$ws.Range("S7").Value = 0.01452962069248976
$ws.Range("T7").Value = 0.01452962069248976

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.68975
$ws.Range("H8").Value = 44.06925
$ws.Range("I8").Value = 0.2692481784748309
$ws.Range("J8").Value = 0.2692481784748309
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.802173666666667
$ws.Range("N8").Value = 5.406521
$ws.Range("O8").Value = 0.03509934175535083
$ws.Range("P8").Value = 0.03509934175535083
$ws.Range("Q8").Value = 26.47348061991666
$ws.Range("R8").Value = 238.26132557925
$ws.Range("S8").Value = 0.009450433833293784
$ws.Range("T8").Value = 0.009450433833293784

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.68975
$ws.Range("H9").Value = 44.06925
$ws.Range("I9").Value = 0.2692481784748309
$ws.Range("J9").Value = 0.2692481784748309
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.7122493333333333
$ws.Range("N9").Value = 2.136748
$ws.Range("O9").Value = 0.01387184999319569
$ws.Range("P9").Value = 0.01387184999319569
$ws.Range("Q9").Value = 10.46276464433333
$ws.Range("R9").Value = 94.16488179899999
$ws.Range("S9").Value = 0.003734970342744036
$ws.Range("T9").Value = 0.003734970342744036

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.967860666666667
$ws.Range("H10").Value = 5.903582
$ws.Range("I10").Value = 0.03606888476606249
$ws.Range("J10").Value = 0.03606888476606249
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.05975733333333
$ws.Range("N10").Value = 138.179272
$ws.Range("O10").Value = 0.8970651351272991
$ws.Range("P10").Value = 0.897065135127299
$ws.Range("Q10").Value = 90.63918477247822
$ws.Range("R10").Value = 815.752662952304
$ws.Range("S10").Value = 0.03235613898655883
$ws.Range("T10").Value = 0.03235613898655882

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.967860666666667
$ws.Range("H11").Value = 5.903582
$ws.Range("I11").Value = 0.03606888476606249
$ws.Range("J11").Value = 0.03606888476606249
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.770761666666667
$ws.Range("N11").Value = 8.312285
$ws.Range("O11").Value = 0.05396367312415441
$ws.Range("P11").Value = 0.0539636731241544
$ws.Range("Q11").Value = 5.452472900541111
$ws.Range("R11").Value = 49.07225610486999
$ws.Range("S11").Value = 0.001946409507468588
$ws.Range("T11").Value = 0.001946409507468588

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.967860666666667
$ws.Range("H12").Value = 5.903582
$ws.Range("I12").Value = 0.03606888476606249
$ws.Range("J12").Value = 0.03606888476606249
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.802173666666667
$ws.Range("N12").Value = 5.406521
$ws.Range("O12").Value = 0.03509934175535083
$ws.Range("P12").Value = 0.03509934175535083
$ws.Range("Q12").Value = 3.546426673135778
$ws.Range("R12").Value = 31.917840058222
$ws.Range("S12").Value = 0.001265994113138395
$ws.Range("T12").Value = 0.001265994113138395

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.967860666666667
$ws.Range("H13").Value = 5.903582
$ws.Range("I13").Value = 0.03606888476606249
$ws.Range("J13").Value = 0.03606888476606249
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.7122493333333333
$ws.Range("N13").Value = 2.136748
$ws.Range("O13").Value = 0.01387184999319569
$ws.Range("P13").Value = 0.01387184999319569
$ws.Range("Q13").Value = 1.401607447926222
$ws.Range("R13").Value = 12.614467031336
$ws.Range("S13").Value = 0.0005003421588966802
$ws.Range("T13").Value = 0.0005003421588966801

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 21.77610566666667
$ws.Range("H14").Value = 65.328317
$ws.Range("I14").Value = 0.3991338712384788
$ws.Range("J14").Value = 0.3991338712384788
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 46.05975733333333
$ws.Range("N14").Value = 138.179272
$ws.Range("O14").Value = 0.8970651351272991
$ws.Range("P14").Value = 0.897065135127299
$ws.Range("Q14").Value = 1003.002142671691
$ws.Range("R14").Value = 9027.019284045224
$ws.Range("S14").Value = 0.3580490801364279
$ws.Range("T14").Value = 0.3580490801364279

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 21.77610566666667
$ws.Range("H15").Value = 65.328317
$ws.Range("I15").Value = 0.3991338712384788
$ws.Range("J15").Value = 0.3991338712384788
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.770761666666667
$ws.Range("N15").Value = 8.312285
$ws.Range("O15").Value = 0.05396367312415441
$ws.Range("P15").Value = 0.0539636731241544
$ws.Range("Q15").Value = 60.33639883048277
$ws.Range("R15").Value = 543.027589474345
$ws.Range("S15").Value = 0.0215387297602916
$ws.Range("T15").Value = 0.0215387297602916

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 21.77610566666667
$ws.Range("H16").Value = 65.328317
$ws.Range("I16").Value = 0.3991338712384788
$ws.Range("J16").Value = 0.3991338712384788
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.802173666666667
$ws.Range("N16").Value = 5.406521
$ws.Range("O16").Value = 0.03509934175535083
$ws.Range("P16").Value = 0.03509934175535083
$ws.Range("Q16").Value = 39.24432419501744
$ws.Range("R16").Value = 353.198917755157
$ws.Range("S16").Value = 0.01400933615273556
$ws.Range("T16").Value = 0.01400933615273556

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 21.77610566666667
$ws.Range("H17").Value = 65.328317
$ws.Range("I17").Value = 0.3991338712384788
$ws.Range("J17").Value = 0.3991338712384788
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.7122493333333333
$ws.Range("N17").Value = 2.136748
$ws.Range("O17").Value = 0.01387184999319569
$ws.Range("P17").Value = 0.01387184999319569
$ws.Range("Q17").Value = 15.51001674367955
$ws.Range("R17").Value = 139.590150693116
$ws.Range("S17").Value = 0.005536725189023663
$ws.Range("T17").Value = 0.005536725189023662

Write-Host "Updated rows 2-17"